$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 283.14285
$ws.Cells.Item(2, 9).Value = 230.5
$ws.Cells.Item(2, 11).Value = 230.5
$ws.Cells.Item(2, 13).Value = -117.5
$ws.Cells.Item(51, 8).Value = 4142.7144
$ws.Cells.Item(51, 10).Value = 4200.4
$ws.Cells.Item(51, 12).Value = 4200.4
$ws.Cells.Item(51, 14).Value = -5168.4
$ws.Cells.Item(64, 8).Value = 8153
$ws.Cells.Item(64, 10).Value = 10833.333
$ws.Cells.Item(64, 12).Value = 10833.333
$ws.Cells.Item(64, 14).Value = -11329.333
$ws.Cells.Item(67, 8).Value = 8153
$ws.Cells.Item(67, 10).Value = 10833.333
$ws.Cells.Item(67, 12).Value = 10833.333
$ws.Cells.Item(67, 14).Value = -12549.333
$ws.Cells.Item(74, 8).Value = 6689.5454
$ws.Cells.Item(74, 9).Value = 5972.778
$ws.Cells.Item(74, 11).Value = 5972.778
$ws.Cells.Item(74, 13).Value = -5036.778
$ws.Cells.Item(77, 8).Value = 6689.5454
$ws.Cells.Item(77, 9).Value = 5972.778
$ws.Cells.Item(77, 11).Value = 29863.89
$ws.Cells.Item(77, 13).Value = -25183.89
$ws.Cells.Item(100, 8).Value = 1875.5
$ws.Cells.Item(100, 9).Value = 1455
$ws.Cells.Item(100, 10).Value = 2296
$ws.Cells.Item(100, 11).Value = 1455
$ws.Cells.Item(100, 12).Value = 2296
$ws.Cells.Item(100, 13).Value = -914
$ws.Cells.Item(100, 14).Value = -3378
$ws.Cells.Item(106, 8).Value = 2422.5
$ws.Cells.Item(106, 9).Value = 2437.3333
$ws.Cells.Item(106, 11).Value = 2437.3333
$ws.Cells.Item(106, 13).Value = -1806.3333
$ws.Cells.Item(129, 8).Value = 4630810.5
$ws.Cells.Item(129, 9).Value = 1328.875
$ws.Cells.Item(129, 11).Value = 3986.625
$ws.Cells.Item(129, 13).Value = 1013.375
$ws.Cells.Item(132, 8).Value = 1599.5416
$ws.Cells.Item(132, 9).Value = 1345.5555
$ws.Cells.Item(132, 11).Value = 4036.6665
$ws.Cells.Item(132, 13).Value = -1506.6665
$ws.Cells.Item(138, 8).Value = 1944.4263
$ws.Cells.Item(138, 9).Value = 1397.9062
$ws.Cells.Item(138, 10).Value = 2547.4827
$ws.Cells.Item(138, 11).Value = 4193.7186
$ws.Cells.Item(138, 12).Value = 7642.4481
$ws.Cells.Item(138, 13).Value = 946.2813999999998
$ws.Cells.Item(138, 14).Value = -17922.4481
$ws.Cells.Item(141, 8).Value = 2895.9375
$ws.Cells.Item(141, 9).Value = 2895.9375
$ws.Cells.Item(141, 11).Value = 8687.8125
$ws.Cells.Item(141, 13).Value = -3507.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2715.4348
$ws.Cells.Item(32, 9).Value = 2433.9768
$ws.Cells.Item(32, 11).Value = 2433.9768
$ws.Cells.Item(32, 13).Value = -2146.9768
$ws.Cells.Item(74, 8).Value = 3783.5854
$ws.Cells.Item(74, 9).Value = 2607.6333
$ws.Cells.Item(74, 10).Value = 6990.727
$ws.Cells.Item(74, 11).Value = 2607.6333
$ws.Cells.Item(74, 12).Value = 6990.727
$ws.Cells.Item(74, 13).Value = -1733.6333
$ws.Cells.Item(74, 14).Value = -8738.726999999999
$ws.Cells.Item(77, 8).Value = 3783.5854
$ws.Cells.Item(77, 9).Value = 2607.6333
$ws.Cells.Item(77, 10).Value = 6990.727
$ws.Cells.Item(77, 11).Value = 13038.1665
$ws.Cells.Item(77, 12).Value = 34953.635
$ws.Cells.Item(77, 13).Value = -8670.166499999999
$ws.Cells.Item(77, 14).Value = -43689.635
$ws.Cells.Item(97, 8).Value = 745
$ws.Cells.Item(97, 9).Value = 860
$ws.Cells.Item(97, 11).Value = 860
$ws.Cells.Item(97, 13).Value = -364

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 249.22223
$ws.Cells.Item(22, 9).Value = 249.22223
$ws.Cells.Item(22, 11).Value = 249.22223
$ws.Cells.Item(22, 13).Value = -76.22223
$ws.Cells.Item(86, 8).Value = 1342
$ws.Cells.Item(86, 9).Value = 1342
$ws.Cells.Item(86, 11).Value = 1342
$ws.Cells.Item(86, 13).Value = -219
$ws.Cells.Item(89, 8).Value = 1342
$ws.Cells.Item(89, 9).Value = 1342
$ws.Cells.Item(89, 11).Value = 6710
$ws.Cells.Item(89, 13).Value = -1094
$ws.Cells.Item(94, 8).Value = 1759
$ws.Cells.Item(94, 9).Value = 1640.7894
$ws.Cells.Item(94, 11).Value = 1640.7894
$ws.Cells.Item(94, 13).Value = -1189.7894

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4519.6206
$ws.Cells.Item(31, 9).Value = 3882.6428
$ws.Cells.Item(31, 10).Value = 5114.1333
$ws.Cells.Item(31, 11).Value = 3882.6428
$ws.Cells.Item(31, 12).Value = 5114.1333
$ws.Cells.Item(31, 13).Value = -3587.6428
$ws.Cells.Item(31, 14).Value = -5704.1333
$ws.Cells.Item(34, 8).Value = 4519.6206
$ws.Cells.Item(34, 9).Value = 3882.6428
$ws.Cells.Item(34, 10).Value = 5114.1333
$ws.Cells.Item(34, 11).Value = 3882.6428
$ws.Cells.Item(34, 12).Value = 5114.1333
$ws.Cells.Item(34, 13).Value = -3680.6428
$ws.Cells.Item(34, 14).Value = -5518.1333
$ws.Cells.Item(74, 8).Value = 38329
$ws.Cells.Item(74, 10).Value = 41336.332
$ws.Cells.Item(74, 12).Value = 41336.332
$ws.Cells.Item(74, 14).Value = -43084.332
$ws.Cells.Item(77, 8).Value = 38329
$ws.Cells.Item(77, 10).Value = 41336.332
$ws.Cells.Item(77, 12).Value = 124008.996
$ws.Cells.Item(77, 14).Value = -132744.996
$ws.Cells.Item(111, 8).Value = 54979
$ws.Cells.Item(111, 10).Value = 54979
$ws.Cells.Item(111, 12).Value = 54979
$ws.Cells.Item(111, 14).Value = -63159
$ws.Cells.Item(112, 8).Value = 79989.5
$ws.Cells.Item(112, 10).Value = 79989.5
$ws.Cells.Item(112, 12).Value = 79989.5
$ws.Cells.Item(112, 14).Value = -82943.5
$ws.Cells.Item(134, 8).Value = 7972.8887
$ws.Cells.Item(134, 9).Value = 6964.143
$ws.Cells.Item(134, 10).Value = 11503.5
$ws.Cells.Item(134, 11).Value = 20892.429
$ws.Cells.Item(134, 12).Value = 34510.5
$ws.Cells.Item(134, 13).Value = -18357.429
$ws.Cells.Item(134, 14).Value = -39580.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 550625.9
$ws.Cells.Item(140, 9).Value = 2471.4
$ws.Cells.Item(140, 11).Value = 7414.200000000001
$ws.Cells.Item(140, 13).Value = -2234.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 759000
$ws.Cells.Item(21, 10).Value = 18000
$ws.Cells.Item(21, 12).Value = 18000
$ws.Cells.Item(21, 14).Value = -18346
$ws.Cells.Item(24, 8).Value = 15379.375
$ws.Cells.Item(24, 9).Value = 10000
$ws.Cells.Item(24, 11).Value = 10000
$ws.Cells.Item(24, 13).Value = -9827
$ws.Cells.Item(30, 8).Value = 759000
$ws.Cells.Item(30, 10).Value = 18000
$ws.Cells.Item(30, 12).Value = 18000
$ws.Cells.Item(30, 14).Value = -18210
$ws.Cells.Item(102, 8).Value = 3128.5715
$ws.Cells.Item(102, 9).Value = 2577.7778
$ws.Cells.Item(102, 10).Value = 4120
$ws.Cells.Item(102, 11).Value = 2577.7778
$ws.Cells.Item(102, 12).Value = 4120
$ws.Cells.Item(102, 13).Value = -955.7777999999998
$ws.Cells.Item(102, 14).Value = -7364

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1679.6666
$ws.Cells.Item(22, 9).Value = 1819.5
$ws.Cells.Item(22, 11).Value = 1819.5
$ws.Cells.Item(22, 13).Value = -1524.5
$ws.Cells.Item(27, 8).Value = 1679.6666
$ws.Cells.Item(27, 9).Value = 1819.5
$ws.Cells.Item(27, 11).Value = 1819.5
$ws.Cells.Item(27, 13).Value = -1712.5
$ws.Cells.Item(101, 8).Value = 76333.336
$ws.Cells.Item(101, 10).Value = 76333.336
$ws.Cells.Item(101, 12).Value = 76333.336
$ws.Cells.Item(101, 14).Value = -82823.336
$ws.Cells.Item(136, 8).Value = 4259.7856
$ws.Cells.Item(136, 9).Value = 3920.1667
$ws.Cells.Item(136, 10).Value = 6297.5
$ws.Cells.Item(136, 11).Value = 11760.5001
$ws.Cells.Item(136, 12).Value = 18892.5
$ws.Cells.Item(136, 13).Value = -9210.500100000001
$ws.Cells.Item(136, 14).Value = -23992.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2134
$ws.Cells.Item(81, 9).Value = 904.4
$ws.Cells.Item(81, 10).Value = 8282
$ws.Cells.Item(81, 11).Value = 1808.8
$ws.Cells.Item(81, 12).Value = 16564
$ws.Cells.Item(81, 13).Value = -747.8
$ws.Cells.Item(81, 14).Value = -18686
$ws.Cells.Item(84, 8).Value = 2134
$ws.Cells.Item(84, 9).Value = 904.4
$ws.Cells.Item(84, 10).Value = 8282
$ws.Cells.Item(84, 11).Value = 9044
$ws.Cells.Item(84, 12).Value = 82820
$ws.Cells.Item(84, 13).Value = -3740
$ws.Cells.Item(84, 14).Value = -93428
$ws.Cells.Item(133, 8).Value = 189999
$ws.Cells.Item(133, 10).Value = 189999
$ws.Cells.Item(133, 12).Value = 189999
$ws.Cells.Item(133, 14).Value = -200119
